$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: McMaster Order - Additional Aluminum, Threaded Rods, Nuts, Wiring Tubing ---
# Copy the date-format (and general row styling) from row 10 down to row 11 first,
# so the new date cell reuses the existing date style instead of Excel minting a
# brand-new number-format style.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("A11").Value = Get-Date -Year 2018 -Month 4 -Day 16 -Hour 0 -Minute 0 -Second 0
$ws.Range("B11").Value = "McMaster Order"
$ws.Range("C11").Value = "Additional Aluminum, Threaded Rods, Nuts, Wiring Tubing"
$ws.Range("D11").Value = "Mcmaster"
$ws.Range("E11").Value = 57.14

# --- Row 12: Amazon Motor Order - DC Brushed Motor 80T ---
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("A12").Value = Get-Date -Year 2018 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("B12").Value = "Amazon Motor Order"
$ws.Range("C12").Value = "DC Brushed Motor 80T"
$ws.Range("D12").Value = "Amazon"
$ws.Range("E12").Value = 15

# Row 11 grew taller (wrapped long "Additional Aluminum..." text) while row 12 stays default.
$ws.Rows.Item(11).RowHeight = 30

# Selection moved to F12 as the last user action before saving.
$ws.Range("F12").Select()

$wb.Application.CutCopyMode = $false
